$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (AA, AB) -----------------------------------------
# Copy formatting from the last existing header cell (Z1) onto the two new
# header cells, then set their text.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)
$ws.Range("AA1").Value = "Área Priv."
$ws.Range("AB1").Value = "Área Com."

# --- New data columns for data row 2 --------------------------------------
$ws.Range("Y2").Copy()
$ws.Range("AA2:AB2").PasteSpecial(-4122)
# "420" looks numeric, so Excel would normally coerce it to a number; force
# a text number-format before assigning so it is stored as a shared string,
# then restore the original (General) look by re-pasting the formatting.
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "420"
$ws.Range("Y2").Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$ws.Range("AB2").Value = "364,29"

$excel.CutCopyMode = 0
